# FN-3460: fix invalid facility utilisation values -- all need to match as
# same facility id for all rows.
#
# Row 5 (UKEF facility ID 20001371, row "Scone GEF"-equivalent) and Row 6
# had facility-limit / facility-utilisation figures that didn't line up
# with the rest of the rows for the same facility id. Correct them so
# every row for facility 20001371 reports the same facility limit, and fix
# the facility utilisation amount accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Facility limit (column E) for rows 5 and 6 should match the other rows
# for this facility id (600000).
$ws.Range("E5").Value = 600000
$ws.Range("E6").Value = 600000

# Facility utilisation (column G) corrected for rows 5 and 6.
$ws.Range("G5").Value = 3938753.8
$ws.Range("G6").Value = 761579.37

# Reflect the cell range the user inspected while fixing these rows.
$ws.Range("E5:H6").Select() | Out-Null
